# Add new columns I (I0) and J (IF) with header + values, mirroring
# the style used by the existing header row (style index 1 -> bold,
# centered, bordered), which we replicate by copying the style from H1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, thin border) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Re-apply header values in case PasteSpecial touched them
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-20
$values = @{
    2  = @(5, 6)
    3  = @(8, 8)
    4  = @(8, 8)
    5  = @(10, 10)
    6  = @(7, 7)
    7  = @(8, 8)
    8  = @(7, 7)
    9  = @(8, 8)
    10 = @(5, 5)
    11 = @(7, 7)
    12 = @(5, 5)
    13 = @(6, 8)
    14 = @(8, 8)
    15 = @(4, 5)
    16 = @(7, 7)
    17 = @(9, 9)
    18 = @(1, 2)
    19 = @(8, 8)
    20 = @(8, 8)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

$wb.Save()
